$d = $word.ActiveDocument

# Each replacement targets a unique cell text of the form NNN×N=
# Using Find/Execute with MatchWholeWord off but exact text (no wildcards),
# restricted to the full document content, since every old string is unique.
$d.Content.Find.Execute("150×7=", $true, $false, $false, $false, $false, $true, 1, $false, "981×8=", 2) | Out-Null
$d.Content.Find.Execute("715×6=", $true, $false, $false, $false, $false, $true, 1, $false, "353×4=", 2) | Out-Null
$d.Content.Find.Execute("518×7=", $true, $false, $false, $false, $false, $true, 1, $false, "733×8=", 2) | Out-Null
$d.Content.Find.Execute("289×2=", $true, $false, $false, $false, $false, $true, 1, $false, "241×4=", 2) | Out-Null
$d.Content.Find.Execute("456×9=", $true, $false, $false, $false, $false, $true, 1, $false, "280×3=", 2) | Out-Null
$d.Content.Find.Execute("329×7=", $true, $false, $false, $false, $false, $true, 1, $false, "626×8=", 2) | Out-Null
$d.Content.Find.Execute("326×9=", $true, $false, $false, $false, $false, $true, 1, $false, "408×5=", 2) | Out-Null
$d.Content.Find.Execute("481×7=", $true, $false, $false, $false, $false, $true, 1, $false, "224×3=", 2) | Out-Null
$d.Content.Find.Execute("398×9=", $true, $false, $false, $false, $false, $true, 1, $false, "674×7=", 2) | Out-Null
$d.Content.Find.Execute("725×4=", $true, $false, $false, $false, $false, $true, 1, $false, "906×5=", 2) | Out-Null
$d.Content.Find.Execute("168×8=", $true, $false, $false, $false, $false, $true, 1, $false, "658×2=", 2) | Out-Null
$d.Content.Find.Execute("738×2=", $true, $false, $false, $false, $false, $true, 1, $false, "193×4=", 2) | Out-Null
$d.Content.Find.Execute("827×9=", $true, $false, $false, $false, $false, $true, 1, $false, "733×4=", 2) | Out-Null
$d.Content.Find.Execute("698×6=", $true, $false, $false, $false, $false, $true, 1, $false, "512×2=", 2) | Out-Null
$d.Content.Find.Execute("287×6=", $true, $false, $false, $false, $false, $true, 1, $false, "298×9=", 2) | Out-Null
$d.Content.Find.Execute("957×8=", $true, $false, $false, $false, $false, $true, 1, $false, "802×8=", 2) | Out-Null
$d.Content.Find.Execute("265×9=", $true, $false, $false, $false, $false, $true, 1, $false, "678×8=", 2) | Out-Null
$d.Content.Find.Execute("898×5=", $true, $false, $false, $false, $false, $true, 1, $false, "602×6=", 2) | Out-Null
$d.Content.Find.Execute("606×8=", $true, $false, $false, $false, $false, $true, 1, $false, "457×7=", 2) | Out-Null
$d.Content.Find.Execute("911×6=", $true, $false, $false, $false, $false, $true, 1, $false, "221×7=", 2) | Out-Null
$d.Content.Find.Execute("466×6=", $true, $false, $false, $false, $false, $true, 1, $false, "328×8=", 2) | Out-Null
$d.Content.Find.Execute("321×8=", $true, $false, $false, $false, $false, $true, 1, $false, "102×2=", 2) | Out-Null
$d.Content.Find.Execute("498×4=", $true, $false, $false, $false, $false, $true, 1, $false, "642×4=", 2) | Out-Null
$d.Content.Find.Execute("842×7=", $true, $false, $false, $false, $false, $true, 1, $false, "705×9=", 2) | Out-Null
$d.Content.Find.Execute("856×8=", $true, $false, $false, $false, $false, $true, 1, $false, "815×5=", 2) | Out-Null

Write-Host "Replacements applied."
